# CCC19 Derived Variables Spreadsheet - add new derived variable "T10a"
# (quarter_median_dx) right after "T10" in the Table1 listing, matching
# the target commit "Minor fixes and new variables".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The table (Table1) is sorted alphabetically by column A ("Variable #").
# "T10a" sorts immediately after "T10" (row 109) and before "T11"
# (currently row 110), so insert a fresh worksheet row at 110 and shift
# everything below it down by one.
$ws.Rows.Item(110).Insert()

$ws.Cells.Item(110, 1).Value = "T10a"
$ws.Cells.Item(110, 4).Value = "Quarter and year of diagnosis, using the median of the inteval as anchor"
$ws.Cells.Item(110, 2).Value = "quarter_median_dx"
$ws.Cells.Item(110, 3).Value = "Time measurements"

# Grow Table1 so the new row is recognised as part of the table (was
# A1:E118, now A1:E119).
$tbl = $ws.ListObjects.Item(1)
$lastRow = $tbl.Range.Rows.Count + 1
$tbl.Resize($ws.Range("A1:E" + $lastRow))

# Reflect the author's on-screen state at save time.
$ws.Range("B110").Select() | Out-Null
$excel.ActiveWindow.Zoom = 230
